$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 630, shifting the existing data
# (old rows 630-671) down to rows 632-673.
$ws.Rows.Item(630).Resize(2).Insert()

# New row 630: 2026/01/11 (Sun)
$ws.Range("A630").NumberFormat = "@"
$ws.Range("A630").Value = "2026/01/11"
$ws.Range("B630").Value = "日"
$ws.Range("C630").Value = 23
$ws.Range("D630").Value = 152

# New row 631: 2026/01/12 (Mon)
$ws.Range("A631").NumberFormat = "@"
$ws.Range("A631").Value = "2026/01/12"
$ws.Range("B631").Value = "月"
$ws.Range("C631").Value = 1
$ws.Range("D631").Value = 132
